# Auto-generated: apply cryptos price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$text) {
    # Prefix with an apostrophe so Excel treats number-looking strings (e.g. "1.001")
    # as literal text instead of silently coercing them to a number; then restore the
    # default cell style so no stray number-format / quote-prefix styling is left behind.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "22.404.10"
Set-TextCell $ws.Range("E2") "  -0.25%  "
Set-TextCell $ws.Range("D3") "1.563.34"
Set-TextCell $ws.Range("E3") "  -0.56%  "
Set-TextCell $ws.Range("E4") "  -0.06%  "
Set-TextCell $ws.Range("D5") "1.001"
Set-TextCell $ws.Range("E5") "  -0.13%  "
Set-TextCell $ws.Range("D6") "287.06"
Set-TextCell $ws.Range("E6") "  -1.35%  "
Set-TextCell $ws.Range("D7") "0.3662"
Set-TextCell $ws.Range("E7") "  -2.38%  "
Set-TextCell $ws.Range("D8") "49.68"
Set-TextCell $ws.Range("E8") "  -0.50%  "
Set-TextCell $ws.Range("D9") "0.3358"
Set-TextCell $ws.Range("E9") "  -1.53%  "
Set-TextCell $ws.Range("D10") "1.129"
Set-TextCell $ws.Range("E10") "  -1.85%  "
Set-TextCell $ws.Range("D11") "0.07414"
Set-TextCell $ws.Range("E11") "  -2.07%  "
Set-TextCell $ws.Range("E12") "  -0.07%  "
Set-TextCell $ws.Range("D13") "20.93"
Set-TextCell $ws.Range("E13") "  -2.35%  "
Set-TextCell $ws.Range("D14") "5.934"
Set-TextCell $ws.Range("E14") "  -1.48%  "
Set-TextCell $ws.Range("D15") "6.870"
Set-TextCell $ws.Range("E15") "  -1.35%  "
Set-TextCell $ws.Range("D16") "1.562.90"
Set-TextCell $ws.Range("E16") "  -0.53%  "
Set-TextCell $ws.Range("D17") "0.00001100"
Set-TextCell $ws.Range("E17") "  -1.99%  "
Set-TextCell $ws.Range("D18") "89.00"
Set-TextCell $ws.Range("D19") "0.06743"
Set-TextCell $ws.Range("E19") "  -0.25%  "
Set-TextCell $ws.Range("E20") "  -0.15%  "
Set-TextCell $ws.Range("D21") "6.311"
Set-TextCell $ws.Range("E21") "  +0.71%  "
Set-TextCell $ws.Range("D22") "16.07"
Set-TextCell $ws.Range("D23") "11.98"
Set-TextCell $ws.Range("E23") "  -2.08%  "
Set-TextCell $ws.Range("D24") "22.399.21"
Set-TextCell $ws.Range("D25") "2.371"
Set-TextCell $ws.Range("E25") "  +0.75%  "
Set-TextCell $ws.Range("D26") "2.539"
Set-TextCell $ws.Range("E26") "  -2.44%  "
Set-TextCell $ws.Range("D27") "149.52"
Set-TextCell $ws.Range("E27") "  +0.68%  "
Set-TextCell $ws.Range("E28") "  -2.62%  "
Set-TextCell $ws.Range("D29") "5.000"
Set-TextCell $ws.Range("E29") "  +0.10%  "
Set-TextCell $ws.Range("D30") "123.15"
Set-TextCell $ws.Range("E30") "  -2.39%  "
Set-TextCell $ws.Range("D31") "1.738.16"
Set-TextCell $ws.Range("E31") "  -0.55%  "
Set-TextCell $ws.Range("D32") "1.061"
Set-TextCell $ws.Range("E32") "  +1.55%  "
Set-TextCell $ws.Range("E33") "  -1.32%  "
Set-TextCell $ws.Range("D34") "1.990"
Set-TextCell $ws.Range("E34") "  +0.26%  "
Set-TextCell $ws.Range("D35") "9.580"
Set-TextCell $ws.Range("E35") "  -2.82%  "
Set-TextCell $ws.Range("E36") "  -2.21%  "
Set-TextCell $ws.Range("D37") "0.02397"
Set-TextCell $ws.Range("E37") "  -2.88%  "
Set-TextCell $ws.Range("D38") "1.310"
Set-TextCell $ws.Range("E38") "  -4.92%  "
Set-TextCell $ws.Range("D39") "0.2221"
Set-TextCell $ws.Range("E39") "  -3.18%  "
Set-TextCell $ws.Range("D40") "0.06383"
Set-TextCell $ws.Range("E40") "  -2.96%  "
Set-TextCell $ws.Range("D41") "5.346"
Set-TextCell $ws.Range("E41") "  -2.82%  "
Set-TextCell $ws.Range("D42") "11.23"
Set-TextCell $ws.Range("E42") "  -1.37%  "
Set-TextCell $ws.Range("D43") "0.6085"
Set-TextCell $ws.Range("E43") "  -3.62%  "
Set-TextCell $ws.Range("E44") "  +0.17%  "
Set-TextCell $ws.Range("E45") "  -2.73%  "
Set-TextCell $ws.Range("D46") "3.768"
Set-TextCell $ws.Range("E46") "  -1.23%  "
Set-TextCell $ws.Range("D47") "0.5746"
Set-TextCell $ws.Range("E47") "  -2.47%  "
Set-TextCell $ws.Range("D48") "2.016"
Set-TextCell $ws.Range("E48") "  -4.05%  "
Set-TextCell $ws.Range("D49") "124.76"
Set-TextCell $ws.Range("E49") "  -3.73%  "
Set-TextCell $ws.Range("D50") "1.219"
Set-TextCell $ws.Range("E50") "  -0.78%  "
Set-TextCell $ws.Range("D51") "0.07235"
Set-TextCell $ws.Range("E51") "  -1.33%  "
